$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 4
$lastRow = 64

# Read the full data block (Ciudad, Casos totales, Casos activos, Recuperados, Muertes)
$rows = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $row = @{
        Name = $ws.Cells.Item($r, 1).Value2
        B    = $ws.Cells.Item($r, 2).Value2
        C    = $ws.Cells.Item($r, 3).Value2
        D    = $ws.Cells.Item($r, 4).Value2
        E    = $ws.Cells.Item($r, 5).Value2
    }
    $rows += $row
}

# Apply the updated case counts (new data for this refresh)
for ($i = 0; $i -lt $rows.Count; $i++) {
    if ($rows[$i].Name -eq "Asturias") {
        $rows[$i].B = 779
        $rows[$i].C = 35
        $rows[$i].D = 719
        $rows[$i].E = 25
    } elseif ($rows[$i].Name -eq "Murcia") {
        $rows[$i].B = 477
        $rows[$i].C = 4
        $rows[$i].D = 467
        $rows[$i].E = 6
    }
}

# Re-sort descending by "Casos totales" (column B), matching the sheet's ranking order
$sorted = $rows | Sort-Object -Property B -Descending

# Write the resorted block back
$r = $firstRow
foreach ($row in $sorted) {
    $ws.Cells.Item($r, 1).Value = $row.Name
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $r++
}

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 21:46"
